$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the two narrow spacer columns (old D:E) so F..L shift left to D..J
$ws.Range("D:E").Delete()

# 2. Move the "Insurance Deduction / Other Regular Deduction / Total Regular
#    Deductions" block down by two rows (bottom-most first to avoid clobbering).
$ws.Range("H22:I22").Cut($ws.Range("H24:I24"))
$ws.Range("H20:I20").Cut($ws.Range("H22:I22"))
$ws.Range("H18:I18").Cut($ws.Range("H20:I20"))
$ws.Range("H18:I18").Clear()

# 3. Repair formulas whose target cells moved (Cut() here does not follow
#    references automatically the way interactive Excel would).
$ws.Range("I24").Formula = "=SUM(I20,I22)"
$ws.Range("D20").Formula = "=D18*I16+I24"

# 4. "Hourly Wage" label now sits immediately beside "Employee Name"; align left.
$ws.Range("D4").HorizontalAlignment = -4131

# 5. New "Period" label in the top right corner of the header.
$ws.Range("J1").Formula = "=CONCATENATE(""Period: "",TEXT(NOW(), ""M/d/yyyy""))"
$ws.Range("J1").HorizontalAlignment = -4152
$ws.Range("J1").Font.Name = "Segoe UI"
$ws.Range("J1").Font.Size = 11
$ws.Range("J1").Font.ThemeColor = 1
$ws.Range("J1").Font.TintAndShade = 0.249977111117893

# 6. Give the (now unmerged) label column C an explicit width.
$ws.Columns("C").ColumnWidth = 29.7109375

# 7. Keep the selection on the cell that used to be selected before the shift.
$ws.Range("D8").Select()
